# "Ajout du formulaire de contact"
# - Rename Feuil2 / Feuil3 into real test-case sheets, fill them with content
# - Add a brand new "Test CU 4 Gerer les information" sheet
# - Move the active tab from sheet1 to the new last sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Color / style helpers (reproduce the workbook's existing font palette)
# ---------------------------------------------------------------------------
$colBlue  = 12611584   # RGB(0,112,192)   -> "FF0070C0" (used for scenario steps)
$colRed   = 255         # RGB(255,0,0)     -> "FFFF0000" (used for "Ko")
$colGreen = 5287936     # RGB(0,176,80)    -> "FF00B050" (used for "Ok")
$xlVAlignCenter = -4108

function Style-Step($rng) {
    $rng.Font.Color = $colBlue
    $rng.VerticalAlignment = $xlVAlignCenter
}

function Style-Header($rng) {
    $rng.Font.Bold = $true
}

function Style-Ok($rng) {
    $rng.Font.Color = $colGreen
}

function Style-Ko($rng) {
    $rng.Font.Color = $colRed
}

# ---------------------------------------------------------------------------
# Rename the existing placeholder sheets
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Test CU 2 Visiter le site"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "Test CU 3 s'authentifier"

# ---------------------------------------------------------------------------
# Sheet 2: "Test CU 2 Visiter le site"
# ---------------------------------------------------------------------------
$ws2.Range("A1").Value = "Nom du test"
$ws2.Range("B1").Value = "Résultat du test"
$ws2.Range("C1").Value = "Observations"
Style-Header $ws2.Range("A1:C1")

$steps2 = @(
    "1. L'internaute rentre dans son moteur de recherche les mots clé de l'association",
    "2. Le moteur de recherche affiche à l'internaute les résultats référencés de sa recherche sur l'association",
    "3. L'internaute trouve avec les résultats affichés l'association qui l'intéresse",
    "4. L'internaute clique sur le lien de l'association",
    "5. Le moteur de recherche lui dirige vers la page de l'association",
    "6. L'internaute visite le site notamment la page d'accueil où il y a la présentation de l'assocation, les projets, le formulaire de contact"
)
for ($i = 0; $i -lt $steps2.Length; $i++) {
    $r = 2 + $i
    $ws2.Range("A$r").Value = $steps2[$i]
    Style-Step $ws2.Range("A$r")
    if ($r -lt 7) {
        $ws2.Range("B$r").Value = "Ko"
        Style-Ko $ws2.Range("B$r")
        $ws2.Range("C$r").Value = "Fonctionnalité nécessitant un déploiement"
    } else {
        $ws2.Range("B$r").Value = "Ok"
        Style-Ok $ws2.Range("B$r")
    }
}

$ws2.Columns.Item(1).ColumnWidth = 107
$ws2.Columns.Item(2).ColumnWidth = 14.5
$ws2.Columns.Item(3).ColumnWidth = 34.796875
$ws2.Range("A1:C1").Select()

# ---------------------------------------------------------------------------
# Sheet 3: "Test CU 3 s'authentifier"
# ---------------------------------------------------------------------------
$ws3.Range("A1").Value = "Nom du test"
$ws3.Range("B1").Value = "Résultat du test"
$ws3.Range("C1").Value = "Observations"
Style-Header $ws3.Range("A1:C1")

$steps3 = @(
    "1. Le membre clique sur le lien de connexion",
    "2. Le site affiche au membre le formulaire de connexion",
    "3. Le membre saisit dans ce formulaire son nom d'utilisateur et son mot de passe et valide ses choix",
    "4. Le site vérifie les informations saisies par l'internaute",
    "5. Le site informe l'internaute de la véracité des informations saisies",
    "6. Le site affiche un message au membre qu'il est connecté",
    "7. Le site redirige le membre vers la page d'accueil"
)
for ($i = 0; $i -lt $steps3.Length; $i++) {
    $r = 2 + $i
    $ws3.Range("A$r").Value = $steps3[$i]
    Style-Step $ws3.Range("A$r")
    $ws3.Range("B$r").Value = "Ok"
    Style-Ok $ws3.Range("B$r")
}

$ws3.Columns.Item(1).ColumnWidth = 81.3984375
$ws3.Columns.Item(2).ColumnWidth = 14.5
$ws3.Range("A1:C1").Select()

# ---------------------------------------------------------------------------
# Sheet 4 (new): "Test CU 4 Gérer les information"
# ---------------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $last)
$ws4.Name = "Test CU 4 Gérer les information"

$ws4.Range("A1").Value = "Nom du test"
$ws4.Range("B1").Value = "Résultat du test"
$ws4.Range("C1").Value = "Observations"
Style-Header $ws4.Range("A1:C1")

$steps4 = @(
    "1. Le membre clique sur le lien de mise à jour de ses informations",
    "2. Le site affiche au membre le formulaire où il y a ses informations saisies lors de l'inscription",
    "3. Le membre saisit dans ce formulaire les changements qu'il souhaite apporter en les saisissant dans les champs et valide ses choix",
    "4. Le site vérifie les informations saisies par l'internaute",
    "5. Le site informe l'internaute de la véracité des informations saisies",
    "6. Le site affiche un message au membre que ces informations saisies sont bien modifiées",
    "7. Le site redirige le membre vers la page d'accueil"
)
$results4 = @("Ok", "Ok", "Ok", "Ok", "Ko", "Ok", "Ko")
for ($i = 0; $i -lt $steps4.Length; $i++) {
    $r = 2 + $i
    $ws4.Range("A$r").Value = $steps4[$i]
    Style-Step $ws4.Range("A$r")
    $ws4.Range("B$r").Value = $results4[$i]
}

$ws4.Columns.Item(1).ColumnWidth = 108.296875
$ws4.Columns.Item(2).ColumnWidth = 14.5
$ws4.Columns.Item(3).ColumnWidth = 12.59765625
$ws4.Range("B7").Select()

# ---------------------------------------------------------------------------
# Sheet 1: selection moves, it is no longer the displayed tab
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B10").Select()

# Sheet 4 ends up active / the displayed tab, matching the authored edit
$ws4.Activate()
$ws4.Range("B7").Select()
